# Refresh the crypto symbol list: update Price (column D) and Volume(1h)
# (column E) for the changed rows. Values are stored as literal text (not
# numbers/percentages) in the source sheet, so each one is written with a
# leading apostrophe to force Excel to keep it as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'292.43"
$ws.Cells.Item(2, 5).Value = "'-6.95%"
$ws.Cells.Item(3, 4).Value = "'40.44"
$ws.Cells.Item(3, 5).Value = "'-1.36%"
$ws.Cells.Item(4, 4).Value = "'5.037"
$ws.Cells.Item(4, 5).Value = "'-2.30%"
$ws.Cells.Item(5, 4).Value = "'0.07322"
$ws.Cells.Item(5, 5).Value = "'-3.59%"
$ws.Cells.Item(6, 4).Value = "'1.527"
$ws.Cells.Item(6, 5).Value = "'-9.06%"
$ws.Cells.Item(7, 4).Value = "'0.9282"
$ws.Cells.Item(7, 5).Value = "'-0.24%"
$ws.Cells.Item(8, 4).Value = "'2.354"
$ws.Cells.Item(8, 5).Value = "'-2.89%"
$ws.Cells.Item(9, 4).Value = "'0.1176"
$ws.Cells.Item(9, 5).Value = "'-1.92%"
$ws.Cells.Item(10, 5).Value = "'-4.10%"
$ws.Cells.Item(11, 4).Value = "'0.04330"
$ws.Cells.Item(11, 5).Value = "'4.83%"
$ws.Cells.Item(12, 4).Value = "'0.08632"
$ws.Cells.Item(12, 5).Value = "'-4.70%"
$ws.Cells.Item(13, 4).Value = "'0.1055"
$ws.Cells.Item(13, 5).Value = "'0.16%"
$ws.Cells.Item(14, 4).Value = "'0.001274"
$ws.Cells.Item(14, 5).Value = "'-1.32%"
$ws.Cells.Item(15, 4).Value = "'0.005772"
$ws.Cells.Item(15, 5).Value = "'-1.02%"
$ws.Cells.Item(16, 4).Value = "'3.337"
$ws.Cells.Item(16, 5).Value = "'0.20%"
$ws.Cells.Item(17, 4).Value = "'4.292"
$ws.Cells.Item(17, 5).Value = "'-0.92%"
$ws.Cells.Item(18, 4).Value = "'0.3288"
$ws.Cells.Item(18, 5).Value = "'-2.04%"
$ws.Cells.Item(19, 4).Value = "'7.968"
$ws.Cells.Item(19, 5).Value = "'4.82%"
$ws.Cells.Item(20, 4).Value = "'0.1390"
$ws.Cells.Item(20, 5).Value = "'3.62%"
$ws.Cells.Item(21, 4).Value = "'0.2742"
$ws.Cells.Item(21, 5).Value = "'-7.13%"
$ws.Cells.Item(22, 4).Value = "'0.03937"
$ws.Cells.Item(22, 5).Value = "'-2.23%"
$ws.Cells.Item(23, 4).Value = "'0.001261"
$ws.Cells.Item(23, 5).Value = "'-1.38%"
$ws.Cells.Item(24, 4).Value = "'0.003785"
$ws.Cells.Item(24, 5).Value = "'-4.72%"
$ws.Cells.Item(25, 5).Value = "'0.76%"
$ws.Cells.Item(26, 4).Value = "'0.0003724"
$ws.Cells.Item(38, 4).Value = "'0.02286"
$ws.Cells.Item(38, 5).Value = "'-5.44%"
$ws.Cells.Item(39, 4).Value = "'0.05034"
$ws.Cells.Item(39, 5).Value = "'-2.37%"
$ws.Cells.Item(40, 4).Value = "'0.005963"
$ws.Cells.Item(40, 5).Value = "'80.56%"
$ws.Cells.Item(41, 4).Value = "'0.007685"
$ws.Cells.Item(41, 5).Value = "'-0.83%"
$ws.Cells.Item(42, 4).Value = "'0.1287"
$ws.Cells.Item(42, 5).Value = "'-0.88%"
$ws.Cells.Item(43, 4).Value = "'0.007350"
$ws.Cells.Item(43, 5).Value = "'-3.28%"
$ws.Cells.Item(44, 4).Value = "'0.008272"
$ws.Cells.Item(44, 5).Value = "'-3.51%"
$ws.Cells.Item(45, 4).Value = "'0.2916"
$ws.Cells.Item(45, 5).Value = "'-13.98%"
$ws.Cells.Item(46, 4).Value = "'0.00006266"
$ws.Cells.Item(46, 5).Value = "'-4.86%"
$ws.Cells.Item(47, 5).Value = "'-0.01%"
$ws.Cells.Item(48, 4).Value = "'0.03007"
$ws.Cells.Item(48, 5).Value = "'-88.81%"
$ws.Cells.Item(49, 4).Value = "'0.00002101"
$ws.Cells.Item(49, 5).Value = "'-0.01%"
$ws.Cells.Item(50, 4).Value = "'0.0002001"
$ws.Cells.Item(50, 5).Value = "'-0.01%"
